# The paragraph originally reads (as visible text):
#   "A cheezy <SPACE><Le Monde mathematical puzzle (hyperlink)> : (which took ..."
# and must become:
#   "A cheezy : (which took ..."
# i.e. the stray space run and the whole hyperlink run must be removed,
# leaving the trailing " : (which took ..." run (which already starts with
# its own leading space) directly after the "cheezy" spell-check markers.

$d = $word.ActiveDocument

# Step 1: remove the lone space run that sits between "cheezy" and the
# hyperlink. Locate "cheezy " and collapse the found range down to just the
# trailing space character before deleting it, so only that single run is
# touched.
$r1 = $d.Content
$r1.Find.Execute("cheezy ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r1.MoveStart(1, 6) | Out-Null
$r1.Delete()

# Step 2: remove the hyperlink text itself (exact text match keeps the
# deletion confined to the hyperlink's own run, so the now-empty run and its
# enclosing <w:hyperlink> wrapper are dropped entirely).
$r2 = $d.Content
$r2.Find.Execute("Le Monde mathematical puzzle", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r2.Delete()
